# Apply "I0 and IF added" edit: add two new columns (I, J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add labeled headers I1 ("I0") and J1 ("IF"), matching the
#     look (bold, centered, bordered) of the other header cells by copying
#     the formatting from the existing "IP" header cell (H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows: populate I2:J77 with the corresponding I0/IF values.
#     Each entry is (row, I0 value, IF value).
$ijData = @(
    @(2, 7, 7),
    @(3, 8, 8),
    @(4, 7, 8),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 8, 8),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 7, 8),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 8, 8),
    @(19, 9, 10),
    @(20, 10, 10),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 9, 10),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 9, 9),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 9, 9),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 8, 8),
    @(40, 7, 7),
    @(41, 10, 10),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 8, 8),
    @(46, 9, 9),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 9, 9),
    @(53, 8, 8),
    @(54, 7, 7),
    @(55, 7, 7),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 8, 8),
    @(59, 8, 8),
    @(60, 7, 7),
    @(61, 8, 8),
    @(62, 8, 8),
    @(63, 8, 8),
    @(64, 7, 7),
    @(65, 8, 8),
    @(66, 8, 8),
    @(67, 9, 9),
    @(68, 9, 9),
    @(69, 8, 8),
    @(70, 9, 9),
    @(71, 9, 9),
    @(72, 5, 5),
    @(73, 6, 6),
    @(74, 3, 3),
    @(75, 6, 6),
    @(76, 3, 3),
    @(77, 3, 3)
)

foreach ($row in $ijData) {
    $rowIndex = $row[0]
    $i0Value = $row[1]
    $ifValue = $row[2]
    $ws.Cells.Item($rowIndex, 9).Value = $i0Value
    $ws.Cells.Item($rowIndex, 10).Value = $ifValue
}

Write-Output "Added I0 and IF columns (I1:J77)"
